{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the same textual substitutions described by the OOXML diff:\n//   - the date line changes from \"2024-12-09 Monday\" to \"2024-12-10 Tuesday\"\n//   - each \"A\u00d7B=C\" answer cell is replaced with a new \"A\u00d7B=C\" value\n//\n// Every old string is unique in the document, so a simple find & replace\n// (via Range.search + Range.insertText(\"Replace\")) for each pair is safe\n// and unambiguous.\n\nconst replacements = [\n  [\"2024-12-09 Monday\", \"2024-12-10 Tuesday\"],\n  [\"705\u00d78=5640\", \"548\u00d72=1096\"],\n  [\"524\u00d77=3668\", \"815\u00d76=4890\"],\n  [\"936\u00d74=3744\", \"983\u00d79=8847\"],\n  [\"392\u00d77=2744\", \"930\u00d76=5580\"],\n  [\"440\u00d73=1320\", \"307\u00d73=921\"],\n  [\"944\u00d73=2832\", \"203\u00d78=1624\"],\n  [\"744\u00d78=5952\", \"261\u00d78=2088\"],\n  [\"534\u00d79=4806\", \"126\u00d77=882\"],\n  [\"201\u00d78=1608\", \"982\u00d76=5892\"],\n  [\"608\u00d75=3040\", \"686\u00d72=1372\"],\n  [\"930\u00d74=3720\", \"279\u00d77=1953\"],\n  [\"618\u00d78=4944\", \"941\u00d76=5646\"],\n  [\"638\u00d78=5104\", \"376\u00d75=1880\"],\n  [\"556\u00d76=3336\", \"718\u00d79=6462\"],\n  [\"966\u00d77=6762\", \"932\u00d78=7456\"],\n  [\"715\u00d77=5005\", \"933\u00d78=7464\"],\n  [\"343\u00d74=1372\", \"403\u00d77=2821\"],\n  [\"878\u00d73=2634\", \"578\u00d79=5202\"],\n  [\"824\u00d72=1648\", \"788\u00d74=3152\"],\n  [\"713\u00d75=3565\", \"248\u00d75=1240\"],\n  [\"963\u00d79=8667\", \"314\u00d76=1884\"],\n  [\"802\u00d78=6416\", \"344\u00d79=3096\"],\n  [\"360\u00d75=1800\", \"200\u00d74=800\"],\n  [\"189\u00d79=1701\", \"376\u00d73=1128\"],\n  [\"711\u00d72=1422\", \"312\u00d73=936\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the same textual substitutions described by the OOXML diff:\n#   - the date line changes from \"2024-12-09 Monday\" to \"2024-12-10 Tuesday\"\n#   - each \"A\u00d7B=C\" answer cell is replaced with a new \"A\u00d7B=C\" value\n#\n# Every old string is unique in the document, so a document-wide\n# Find/Replace (wdReplaceAll, case-sensitive) for each pair is safe and\n# unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-12-09 Monday\", \"2024-12-10 Tuesday\"),\n    @(\"705\u00d78=5640\", \"548\u00d72=1096\"),\n    @(\"524\u00d77=3668\", \"815\u00d76=4890\"),\n    @(\"936\u00d74=3744\", \"983\u00d79=8847\"),\n    @(\"392\u00d77=2744\", \"930\u00d76=5580\"),\n    @(\"440\u00d73=1320\", \"307\u00d73=921\"),\n    @(\"944\u00d73=2832\", \"203\u00d78=1624\"),\n    @(\"744\u00d78=5952\", \"261\u00d78=2088\"),\n    @(\"534\u00d79=4806\", \"126\u00d77=882\"),\n    @(\"201\u00d78=1608\", \"982\u00d76=5892\"),\n    @(\"608\u00d75=3040\", \"686\u00d72=1372\"),\n    @(\"930\u00d74=3720\", \"279\u00d77=1953\"),\n    @(\"618\u00d78=4944\", \"941\u00d76=5646\"),\n    @(\"638\u00d78=5104\", \"376\u00d75=1880\"),\n    @(\"556\u00d76=3336\", \"718\u00d79=6462\"),\n    @(\"966\u00d77=6762\", \"932\u00d78=7456\"),\n    @(\"715\u00d77=5005\", \"933\u00d78=7464\"),\n    @(\"343\u00d74=1372\", \"403\u00d77=2821\"),\n    @(\"878\u00d73=2634\", \"578\u00d79=5202\"),\n    @(\"824\u00d72=1648\", \"788\u00d74=3152\"),\n    @(\"713\u00d75=3565\", \"248\u00d75=1240\"),\n    @(\"963\u00d79=8667\", \"314\u00d76=1884\"),\n    @(\"802\u00d78=6416\", \"344\u00d79=3096\"),\n    @(\"360\u00d75=1800\", \"200\u00d74=800\"),\n    @(\"189\u00d79=1701\", \"376\u00d73=1128\"),\n    @(\"711\u00d72=1422\", \"312\u00d73=936\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$old, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$new, [ref]2)\n}\n"}
